# "Generate Report for Archive"
# - Update the localization status text from "Ready for handoff" to "In Translation"
#   on every sheet that reports it (Overview!E2/F2, zh-cn!C2, de-de!C2).
# - Narrow the now-shorter "Status" columns to match the refreshed content width.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Overview")
$ws2 = $wb.Worksheets.Item("zh-cn")
$ws3 = $wb.Worksheets.Item("de-de")

$newStatus = "In Translation"

$ws1.Range("E2").Value = $newStatus
$ws1.Range("F2").Value = $newStatus
$ws2.Range("C2").Value = $newStatus
$ws3.Range("C2").Value = $newStatus

# Resize the Status columns to their refreshed (narrower) content width.
$newColumnWidth = 12.5

$ws1.Columns.Item(5).ColumnWidth = $newColumnWidth
$ws1.Columns.Item(6).ColumnWidth = $newColumnWidth
$ws2.Columns.Item(3).ColumnWidth = $newColumnWidth
$ws3.Columns.Item(3).ColumnWidth = $newColumnWidth
